$d = $word.ActiveDocument

# 1. Change "4" to "5" in "Assignment 24.1" -> "Assignment 25.1"
$d.Content.Find.Execute("Assignment 24.1", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Assignment 25.1", 2)

Write-Host "Done"
